$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-29 down to 10-30
$ws.Rows(9).Insert()

# Populate the new row 9 with the new weekly data entry
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44525
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 100112022
$ws.Range("G9").Value = "Arveja Verde"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 12500
$ws.Range("N9").Value = "`$/saco 25 kilos"
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 500
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
